$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (2020-08-19 .. 2020-09-01) appended below the existing data,
# continuing the DateTime / Scheduled flights / Tracked flights / percent table.
$data = @(
    @("2020-08-19", 50, 48),
    @("2020-08-20", 52, 48),
    @("2020-08-21", 58, 52),
    @("2020-08-22", 51, 48),
    @("2020-08-23", 48, 44),
    @("2020-08-24", 56, 52),
    @("2020-08-25", 45, 42),
    @("2020-08-26", 52, 46),
    @("2020-08-27", 71, 64),
    @("2020-08-28", 55, 50),
    @("2020-08-29", 46, 43),
    @("2020-08-30", 51, 46),
    @("2020-08-31", 52, 49),
    @("2020-09-01", 45, 44)
)

$lastRow = 135
$firstNewRow = $lastRow + 1
$rowCount = $data.Count
$lastNewRow = $lastRow + $rowCount

# Carry the formatting (text number format on A, integer on B/C, percent formula
# style on D) down from the last existing row onto the newly inserted block
# before writing any values, so the new cells pick up styles 1/2/2/3 instead of
# Excel's default "General" style - and so typed date-looking text stays text.
$ws.Range("A" + $lastRow + ":D" + $lastRow).Copy()
$ws.Range("A" + $firstNewRow + ":D" + $lastNewRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$row = $firstNewRow
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Range("D" + $row).Formula = "=C" + $row + "/B" + $row
    $row = $row + 1
}

# Match the author's final view state: scrolled down with G144 selected.
$ws.Range("G144").Select() | Out-Null
